$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "Data and Statistics"

# Header row
$ws.Range("A1").Value = "Carrier"
$ws.Range("B1").Value = "Time Flight"

# Carrier rows (reordered + times reformatted)
$ws.Range("A2").Value = "TK"
$ws.Range("B2").Value = "05:50"

$ws.Range("A3").Value = "SU"
$ws.Range("B3").Value = "06:00"

$ws.Range("A4").Value = "S7"
$ws.Range("B4").Value = "06:30"

$ws.Range("A5").Value = "BA"
$ws.Range("B5").Value = "08:05"

# Clear old column C data from rows 7/8 (was used before, no longer needed)
$ws.Range("C7").ClearContents()
$ws.Range("C8").ClearContents()

# Row 7 left blank as a separator
$ws.Range("A7").ClearContents()
$ws.Range("B7").ClearContents()

# Statistics block, now starting at row 8, two columns wide
$ws.Range("A8").Value = "Statistic"
$ws.Range("B8").Value = "Value"

$ws.Range("A9").Value = "Difference"
$ws.Range("B9").Value = "650.00"

$ws.Range("A10").Value = "Median"
$ws.Range("B10").Value = "14200.00"

$ws.Range("A11").Value = "Average"
$ws.Range("B11").Value = "13550.00"
